$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
  @('ECs', 'Efnb3', 'Ephb3', 'ECs', 1, 0.3333333333333333, 0.1533166666666667, 0.45995, 0.1101680340964018, 0.1101680340964017, 3, 1, 0.7925460000000001, 2.377638, 0.02896056603749454, 0.02896056603749454, 0.1215105109, 1.0935945981, 0.003190528626669794, 0.003190528626669793),
  @('ECs', 'Efnb3', 'Ephb3', 'FAPs', 1, 0.3333333333333333, 0.1533166666666667, 0.45995, 0.1101680340964018, 0.1101680340964017, 3, 1, 24.90258766666667, 74.707763, 0.9099699381802406, 0.9099699381802406, 3.817981732427778, 34.36183559185, 0.1002495991761413, 0.1002495991761413),
  @('ECs', 'Efnb3', 'Ephb3', 'Inflammatory-Mac', 1, 0.3333333333333333, 0.1533166666666667, 0.45995, 0.1101680340964018, 0.1101680340964017, 2, 0.6666666666666666, 0.100159, 0.300477, 0.003659928046762479, 0.003659928046762478, 0.01535604401666667, 0.13820439615, 0.0004032070778461059, 0.0004032070778461057),
  @('ECs', 'Efnb3', 'Ephb3', 'MuSCs', 1, 0.3333333333333333, 0.1533166666666667, 0.45995, 0.1101680340964018, 0.1101680340964017, 3, 1, 1.537356666666667, 4.61207, 0.05617682666770443, 0.05617682666770443, 0.2357023996111111, 2.1213215965, 0.006188890555755313, 0.006188890555755313),
  @('ECs', 'Efnb3', 'Ephb3', 'Resolving-Mac', 1, 0.3333333333333333, 0.1533166666666667, 0.45995, 0.1101680340964018, 0.1101680340964017, 2, 0.6666666666666666, 0.03373566666666667, 0.101207, 0.001232741067797836, 0.001232741067797835, 0.005172239961111112, 0.04655015965000001, 0.0001358086599891867, 0.0001358086599891866),
  @('FAPs', 'Efnb3', 'Ephb3', 'ECs', 3, 1, 0.9559960000000002, 2.867988, 0.6869455370628789, 0.6869455370628788, 3, 1, 0.7925460000000001, 2.377638, 0.02896056603749454, 0.02896056603749454, 0.7576708058160002, 6.819037252344001, 0.01989433159027166, 0.01989433159027165),
  @('FAPs', 'Efnb3', 'Ephb3', 'FAPs', 3, 1, 0.9559960000000002, 2.867988, 0.6869455370628789, 0.6869455370628788, 3, 1, 24.90258766666667, 74.707763, 0.9099699381802406, 0.9099699381802406, 23.80677419898267, 214.260967790844, 0.6250997878943001, 0.6250997878943),
  @('FAPs', 'Efnb3', 'Ephb3', 'Inflammatory-Mac', 3, 1, 0.9559960000000002, 2.867988, 0.6869455370628789, 0.6869455370628788, 2, 0.6666666666666666, 0.100159, 0.300477, 0.003659928046762479, 0.003659928046762478, 0.09575160336400002, 0.8617644302760001, 0.002514171237694744, 0.002514171237694744),
  @('FAPs', 'Efnb3', 'Ephb3', 'MuSCs', 3, 1, 0.9559960000000002, 2.867988, 0.6869455370628789, 0.6869455370628788, 3, 1, 1.537356666666667, 4.61207, 0.05617682666770443, 0.05617682666770443, 1.469706823906667, 13.22736141516, 0.03859042036573448, 0.03859042036573447),
  @('FAPs', 'Efnb3', 'Ephb3', 'Resolving-Mac', 3, 1, 0.9559960000000002, 2.867988, 0.6869455370628789, 0.6869455370628788, 2, 0.6666666666666666, 0.03373566666666667, 0.101207, 0.001232741067797836, 0.001232741067797835, 0.03225116239066668, 0.290260461516, 0.0008468259748778509, 0.0008468259748778506),
  @('MuSCs', 'Efnb3', 'Ephb3', 'ECs', 2, 0.6666666666666666, 0.1265133333333333, 0.37954, 0.09090808927263468, 0.09090808927263468, 3, 1, 0.7925460000000001, 2.377638, 0.02896056603749454, 0.02896056603749454, 0.10026763628, 0.90240872652, 0.002632749722722586, 0.002632749722722586),
  @('MuSCs', 'Efnb3', 'Ephb3', 'FAPs', 2, 0.6666666666666666, 0.1265133333333333, 0.37954, 0.09090808927263468, 0.09090808927263468, 3, 1, 24.90258766666667, 74.707763, 0.9099699381802406, 0.9099699381802406, 3.150509374335555, 28.35458436902, 0.08272362837550318, 0.08272362837550318),
  @('MuSCs', 'Efnb3', 'Ephb3', 'Inflammatory-Mac', 2, 0.6666666666666666, 0.1265133333333333, 0.37954, 0.09090808927263468, 0.09090808927263468, 2, 0.6666666666666666, 0.100159, 0.300477, 0.003659928046762479, 0.003659928046762478, 0.01267144895333333, 0.11404304058, 0.0003327170656065029, 0.0003327170656065028),
  @('MuSCs', 'Efnb3', 'Ephb3', 'MuSCs', 2, 0.6666666666666666, 0.1265133333333333, 0.37954, 0.09090808927263468, 0.09090808927263468, 3, 1, 1.537356666666667, 4.61207, 0.05617682666770443, 0.05617682666770443, 0.1944961164222222, 1.7504650478, 0.005106927973760999, 0.005106927973760999),
  @('MuSCs', 'Efnb3', 'Ephb3', 'Resolving-Mac', 2, 0.6666666666666666, 0.1265133333333333, 0.37954, 0.09090808927263468, 0.09090808927263468, 2, 0.6666666666666666, 0.03373566666666667, 0.101207, 0.001232741067797836, 0.001232741067797835, 0.004268011642222223, 0.03841210478, 0.0001120661350414086, 0.0001120661350414086),
  @('Resolving-Mac', 'Efnb3', 'Ephb3', 'ECs', 1, 0.3333333333333333, 0.155836, 0.467508, 0.1119783395680848, 0.1119783395680848, 3, 1, 0.7925460000000001, 2.377638, 0.02896056603749454, 0.02896056603749454, 0.123507198456, 1.111564786104, 0.003242956097830507, 0.003242956097830506),
  @('Resolving-Mac', 'Efnb3', 'Ephb3', 'FAPs', 1, 0.3333333333333333, 0.155836, 0.467508, 0.1119783395680848, 0.1119783395680848, 3, 1, 24.90258766666667, 74.707763, 0.9099699381802406, 0.9099699381802406, 3.880719651622667, 34.926476864604, 0.1018969227342961, 0.1018969227342961),
  @('Resolving-Mac', 'Efnb3', 'Ephb3', 'Inflammatory-Mac', 1, 0.3333333333333333, 0.155836, 0.467508, 0.1119783395680848, 0.1119783395680848, 2, 0.6666666666666666, 0.100159, 0.300477, 0.003659928046762479, 0.003659928046762478, 0.015608377924, 0.140475401316, 0.0004098326656151261, 0.000409832665615126),
  @('Resolving-Mac', 'Efnb3', 'Ephb3', 'MuSCs', 1, 0.3333333333333333, 0.155836, 0.467508, 0.1119783395680848, 0.1119783395680848, 3, 1, 1.537356666666667, 4.61207, 0.05617682666770443, 0.05617682666770443, 0.2395755135066667, 2.15617962156, 0.006290587772453647, 0.006290587772453646),
  @('Resolving-Mac', 'Efnb3', 'Ephb3', 'Resolving-Mac', 1, 0.3333333333333333, 0.155836, 0.467508, 0.1119783395680848, 0.1119783395680848, 2, 0.6666666666666666, 0.03373566666666667, 0.101207, 0.001232741067797836, 0.001232741067797835, 0.005257231350666668, 0.047315082156, 0.0001380402978893895, 0.0001380402978893894)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $rowNum = $i + 2
    $rowData = $data[$i]
    for ($j = 0; $j -lt $rowData.Length; $j++) {
        $ws.Cells.Item($rowNum, $j + 1).Value = $rowData[$j]
    }
}

$ws.Range("A1").Value = "Sending cluster"